# Insert a new data row at row 164 (pushing existing rows 164-183 down to
# 165-184) on the "Mango" price sheet. The new row carries the same static
# category columns (A,B,C,E..L,N..T) as its neighbours, with a fresh date
# (45142) and volume (300).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 164:183 down by one row, creating a blank row 164.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164.
$ws.Cells.Item(164, 1).Value  = 5
$ws.Cells.Item(164, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(164, 3).Value  = "Maule"
$ws.Cells.Item(164, 4).Value  = 45142
$ws.Cells.Item(164, 5).Value  = 7
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100108
$ws.Cells.Item(164, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value  = 100108002
$ws.Cells.Item(164, 10).Value = "Mango"
$ws.Cells.Item(164, 11).Value = "Sin especificar"
$ws.Cells.Item(164, 12).Value = "Primera"
$ws.Cells.Item(164, 13).Value = 300
$ws.Cells.Item(164, 14).Value = 8000
$ws.Cells.Item(164, 15).Value = 8000
$ws.Cells.Item(164, 16).Value = 8000
$ws.Cells.Item(164, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(164, 18).Value = "Brasil"
$ws.Cells.Item(164, 19).Value = 2000
$ws.Cells.Item(164, 20).Value = 4
